# Update Titan_Profits leve-profit computed columns (H:N) across sheets
# per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 172.92857
$ws.Range("I39").Value = 50.363636
$ws.Range("J39").Value = 622.3333
$ws.Range("K39").Value = 151.090908
$ws.Range("L39").Value = 1866.9999
$ws.Range("M39").Value = 144.909092
$ws.Range("N39").Value = -2458.9999

$ws.Range("H43").Value = 1190.3636
$ws.Range("J43").Value = 1100.6666
$ws.Range("L43").Value = 1100.6666
$ws.Range("N43").Value = -1238.6666

$ws.Range("H70").Value = 2150.3333
$ws.Range("I70").Value = 1462.375
$ws.Range("J70").Value = 2936.5715
$ws.Range("K70").Value = 4387.125
$ws.Range("L70").Value = 8809.7145
$ws.Range("M70").Value = -4117.125
$ws.Range("N70").Value = -9349.7145

$ws.Range("H73").Value = 2150.3333
$ws.Range("I73").Value = 1462.375
$ws.Range("J73").Value = 2936.5715
$ws.Range("K73").Value = 4387.125
$ws.Range("L73").Value = 8809.7145
$ws.Range("M73").Value = -3451.125
$ws.Range("N73").Value = -10681.7145

$ws.Range("H88").Value = 2694.4666
$ws.Range("I88").Value = 1524.75
$ws.Range("K88").Value = 1524.75
$ws.Range("M88").Value = -1118.75

$ws.Range("H91").Value = 2694.4666
$ws.Range("I91").Value = 1524.75
$ws.Range("K91").Value = 1524.75
$ws.Range("M91").Value = -120.75

$ws.Range("H98").Value = 704175.6
$ws.Range("I98").Value = 1124791.2
$ws.Range("J98").Value = 3149.6667
$ws.Range("K98").Value = 1124791.2
$ws.Range("L98").Value = 3149.6667
$ws.Range("M98").Value = -1123293.2
$ws.Range("N98").Value = -6145.6667

$ws.Range("H122").Value = 704175.6
$ws.Range("I122").Value = 1124791.2
$ws.Range("J122").Value = 3149.6667
$ws.Range("K122").Value = 3374373.6
$ws.Range("L122").Value = 9449.000100000001
$ws.Range("M122").Value = -3371923.6
$ws.Range("N122").Value = -14349.0001

$ws.Range("H137").Value = 83335330
$ws.Range("I137").Value = 111112664
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 333337992
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = -333335442
$ws.Range("N137").Value = -15099.9999

$ws.Range("H138").Value = 8460524
$ws.Range("I138").Value = 2022519.8
$ws.Range("J138").Value = 14709175
$ws.Range("K138").Value = 6067559.4
$ws.Range("L138").Value = 44127525
$ws.Range("M138").Value = -6062419.4
$ws.Range("N138").Value = -44137805

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2464.9473
$ws.Range("I61").Value = 2233.3845
$ws.Range("K61").Value = 2233.3845
$ws.Range("M61").Value = -2021.3845

$ws.Range("H74").Value = 3674.851
$ws.Range("I74").Value = 1039.9429
$ws.Range("J74").Value = 11360
$ws.Range("K74").Value = 1039.9429
$ws.Range("L74").Value = 11360
$ws.Range("M74").Value = -165.9429
$ws.Range("N74").Value = -13108

$ws.Range("H77").Value = 3674.851
$ws.Range("I77").Value = 1039.9429
$ws.Range("J77").Value = 11360
$ws.Range("K77").Value = 5199.7145
$ws.Range("L77").Value = 56800
$ws.Range("M77").Value = -831.7145
$ws.Range("N77").Value = -65536

$ws.Range("H122").Value = 1669.5
$ws.Range("I122").Value = 1372.3334
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 4117.0002
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -1667.0002
$ws.Range("N122").Value = -10799.9998

$ws.Range("H136").Value = 2464.9473
$ws.Range("I136").Value = 2233.3845
$ws.Range("K136").Value = 6700.1535
$ws.Range("M136").Value = -4150.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1324.8334
$ws.Range("I107").Value = 1332.5385
$ws.Range("J107").Value = 1304.8
$ws.Range("K107").Value = 1332.5385
$ws.Range("L107").Value = 1304.8
$ws.Range("M107").Value = 587.4614999999999
$ws.Range("N107").Value = -5144.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3610.38
$ws.Range("I31").Value = 1578.0344
$ws.Range("K31").Value = 1578.0344
$ws.Range("M31").Value = -1283.0344

$ws.Range("H34").Value = 3610.38
$ws.Range("I34").Value = 1578.0344
$ws.Range("K34").Value = 1578.0344
$ws.Range("M34").Value = -1376.0344

$ws.Range("H122").Value = 1852.3
$ws.Range("I122").Value = 1183.9474
$ws.Range("J122").Value = 3006.7273
$ws.Range("K122").Value = 3551.8422
$ws.Range("L122").Value = 9020.1819
$ws.Range("M122").Value = -1101.8422
$ws.Range("N122").Value = -13920.1819

$ws.Range("H132").Value = 4903913.5
$ws.Range("I132").Value = 6667944.5
$ws.Range("J132").Value = 3827.889
$ws.Range("K132").Value = 20003833.5
$ws.Range("L132").Value = 11483.667
$ws.Range("M132").Value = -20001303.5
$ws.Range("N132").Value = -16543.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2879.907
$ws.Range("I102").Value = 1758.6666
$ws.Range("J102").Value = 6580
$ws.Range("K102").Value = 1758.6666
$ws.Range("L102").Value = 6580
$ws.Range("M102").Value = -136.6666
$ws.Range("N102").Value = -9824

$ws.Range("H111").Value = 31950
$ws.Range("J111").Value = 31950
$ws.Range("L111").Value = 31950
$ws.Range("N111").Value = -38084

$ws.Range("H126").Value = 2845.8
$ws.Range("I126").Value = 2830
$ws.Range("J126").Value = 2849.75
$ws.Range("K126").Value = 8490
$ws.Range("L126").Value = 8549.25
$ws.Range("M126").Value = -6020
$ws.Range("N126").Value = -13489.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3865.8262
$ws.Range("I122").Value = 3364.6667
$ws.Range("J122").Value = 3941
$ws.Range("K122").Value = 10094.0001
$ws.Range("L122").Value = 11823
$ws.Range("M122").Value = -7644.000100000001
$ws.Range("N122").Value = -16723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 60002
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 60002
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 60002
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -60226

$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50812

$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52808

$ws.Range("H100").Value = 1167.3334
$ws.Range("I100").Value = 1167.3334
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2334.6668
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1793.6668
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 2243.0715
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2567.1667
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 7701.500100000001
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -12601.5001

$ws.Range("H126").Value = 72338.71000000001
$ws.Range("I126").Value = 125349.125
$ws.Range("J126").Value = 1658.1666
$ws.Range("K126").Value = 376047.375
$ws.Range("L126").Value = 4974.4998
$ws.Range("M126").Value = -373577.375
$ws.Range("N126").Value = -9914.4998

